# Temporarily limit availability of forestry residues to 0
# - Insert a new "AGR" worksheet right after "Regions" (before "PWR").
# - It carries a ~TFM_INS block that sets ACT_BND on ABIOFRSR* to 0 (2018) then 5 afterwards.
# - The previously active sheet (PWR) keeps a new remembered selection, and AGR becomes
#   the active tab.

$wb = $excel.ActiveWorkbook

$regions = $wb.Worksheets.Item("Regions")

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $regions)
$ws.Name = "AGR"

# --- Row 2: section header ------------------------------------------------
$title = $ws.Range("B2")
$title.Value = "~TFM_INS"
$title.Font.Bold = $true
$title.Font.Name = "Arial"
$title.Font.Size = 10
$title.Font.ColorIndex = 12

# --- Row 3: column headers -------------------------------------------------
$beHdr = $ws.Range("B3:E3")
$beHdr.Font.Bold = $true
$beHdr.Font.Name = "Arial"
$beHdr.Font.Size = 10
$beHdr.Interior.ColorIndex = 43
$beHdr.Borders.Item(8).LineStyle = 1
$beHdr.Borders.Item(8).Weight = 2
$beHdr.Borders.Item(9).LineStyle = 1
$beHdr.Borders.Item(9).Weight = -4138

$ws.Range("B3").Value = "TimeSlice"
$ws.Range("C3").Value = "LimType"
$ws.Range("D3").Value = "Attribute"
$ws.Range("E3").Value = "Year"

$fgHdr = $ws.Range("F3:G3")
$fgHdr.Font.Bold = $true
$fgHdr.Font.Name = "Arial"
$fgHdr.Font.Size = 10
$fgHdr.Interior.ColorIndex = 44
$fgHdr.Borders.Item(8).LineStyle = 1
$fgHdr.Borders.Item(8).Weight = 2
$fgHdr.Borders.Item(9).LineStyle = 1
$fgHdr.Borders.Item(9).Weight = -4138

$ws.Range("F3").Formula = "=IF(Regions!C`$3<>`"`",Regions!C`$3,`"*`")"
$ws.Range("G3").Formula = "=IF(Regions!D`$3<>`"`",Regions!D`$3,`"*`")"

$hjHdr = $ws.Range("H3:J3")
$hjHdr.Font.Name = "Arial"
$hjHdr.Font.Size = 10
$hjHdr.Interior.ColorIndex = 43
$hjHdr.Borders.Item(8).LineStyle = 1
$hjHdr.Borders.Item(8).Weight = 2
$hjHdr.Borders.Item(9).LineStyle = 1
$hjHdr.Borders.Item(9).Weight = -4138

$ws.Range("H3").Value = "UC_Desc"
$ws.Range("I3").Value = "Pset_PN"
$ws.Range("J3").Value = "Cset_CN"

# --- Rows 4-5: data ----------------------------------------------------
$ws.Range("D4").Value = "ACT_BND"
$ws.Range("E4").Value = 2018
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("I4").Value = "ABIOFRSR*"

$ws.Range("D5").Value = "ACT_BND"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 5
$ws.Range("I5").Value = "ABIOFRSR*"

# --- Restore PWR's remembered selection, then make AGR the active tab -----
$pwr = $wb.Worksheets.Item("PWR")
$pwr.Activate()
[void]$pwr.Range("B7:J10").Select()

[void]$ws.Activate()
[void]$ws.Range("E6").Select()
